# Automatic update of files.
#
# Rows 2-4 hold species-observation records keyed by Id (col A). The source
# sync re-ordered these three records: the record that was in row 3 is now
# in row 2, the one in row 4 is now in row 3, and the one in row 2 is now
# in row 4 (columns A, B, E, F, G, H, Q, R carry the per-record data; all
# other columns/rows are identical across the three records and stay put).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 <- old row 3's record
$ws.Cells.Item(2, 1).Value  = 111525236          # A2 Id
$ws.Cells.Item(2, 2).Value  = 77515              # B2 Taxonsorteringsordning
$ws.Cells.Item(2, 5).Value  = 6425               # E2 TaxonId
$ws.Cells.Item(2, 6).Value  = "Garnlav"          # F2 Artnamn
$ws.Cells.Item(2, 7).Value  = "Alectoria sarmentosa"  # G2 Vetenskapligt namn
$ws.Cells.Item(2, 8).Value  = "(Ach.) Ach."      # H2 Auktor
$ws.Cells.Item(2, 17).Value = 404471.921986955   # Q2 Ost
$ws.Cells.Item(2, 18).Value = 6706721.507764764  # R2 Nord

# Row 3 <- old row 4's record
$ws.Cells.Item(3, 1).Value  = 111525234          # A3 Id
$ws.Cells.Item(3, 2).Value  = 78107              # B3 Taxonsorteringsordning
$ws.Cells.Item(3, 5).Value  = 6453               # E3 TaxonId
$ws.Cells.Item(3, 6).Value  = "Vedskivlav"       # F3 Artnamn
$ws.Cells.Item(3, 7).Value  = "Hertelidea botryosa"  # G3 Vetenskapligt namn
$ws.Cells.Item(3, 8).Value  = "(Fr.) Printzen & Kantvilas"  # H3 Auktor
$ws.Cells.Item(3, 17).Value = 404542.4279557943  # Q3 Ost
$ws.Cells.Item(3, 18).Value = 6706754.164833304  # R3 Nord

# Row 4 <- old row 2's record
$ws.Cells.Item(4, 1).Value  = 111525237          # A4 Id
$ws.Cells.Item(4, 2).Value  = 77515              # B4 Taxonsorteringsordning
$ws.Cells.Item(4, 5).Value  = 6425               # E4 TaxonId
$ws.Cells.Item(4, 6).Value  = "Garnlav"          # F4 Artnamn
$ws.Cells.Item(4, 7).Value  = "Alectoria sarmentosa"  # G4 Vetenskapligt namn
$ws.Cells.Item(4, 8).Value  = "(Ach.) Ach."      # H4 Auktor
$ws.Cells.Item(4, 17).Value = 404472.9612160316  # Q4 Ost
$ws.Cells.Item(4, 18).Value = 6706723.452812355  # R4 Nord
